# Add a new "function" column (D) recording which curve type (gaussian /
# lorentzian) was fit for each detected particle, and rewrite the existing
# "centers" (B) / "sigmas" (C) numeric columns as text so every numeric
# value is stored with its full float precision (matches the Python
# str(float) formatting used by wavelength_calculation_series.py).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- header row --------------------------------------------------------
# New "function" header cell should look like its B1/C1 neighbours (bold,
# bordered "header" cell style), so copy their formatting over first.
$ws.Range("C1").Copy()
$ws.Range("D1").PasteSpecial(-4122)   # xlPasteFormats
$ws.Range("D1").Value = "function"

# --- data rows -----------------------------------------------------------
# Row -> (centers, sigmas, function)
$rows = @(
    @{ Row = 2; B = "682.0097710126338";  C = "22.836389535853694"; D = "gaussian" },
    @{ Row = 3; B = "683.2946254117127";  C = "22.310971228927166"; D = "gaussian" },
    @{ Row = 4; B = "681.0045398935124";  C = "23.082143914472518"; D = "lorentzian" },
    @{ Row = 5; B = "682.3968611566785";  C = "22.337690890305094"; D = "gaussian" },
    @{ Row = 6; B = "683.5460627536477";  C = "26.623546419658084"; D = "gaussian" },
    @{ Row = 7; B = "683.6637693095828";  C = "29.287129892903717"; D = "gaussian" },
    @{ Row = 8; B = "683.2815837865073";  C = "29.256428267484743"; D = "gaussian" },
    @{ Row = 9; B = "682.3600036328893";  C = "34.44729360162933";  D = "lorentzian" }
)

foreach ($r in $rows) {
    $bCell = $ws.Cells.Item($r.Row, 2)
    $cCell = $ws.Cells.Item($r.Row, 3)
    $dCell = $ws.Cells.Item($r.Row, 4)

    # Force text storage (instead of numeric) for the center/sigma values by
    # flipping the cell to a text number format before writing the string,
    # then clear the format back off again so no extra cell formatting is
    # left behind on the cell itself.
    $bCell.NumberFormat = "@"
    $bCell.Value = $r.B
    $bCell.ClearFormats()

    $cCell.NumberFormat = "@"
    $cCell.Value = $r.C
    $cCell.ClearFormats()

    $dCell.Value = $r.D
}
